# Users sheet edit:
#  - remove the "Users" header row, closing the blank row-2 gap so the
#    name list starts at row 1 with no gaps
#  - append two more names (Yung, Alex) after the existing list
#  - append four "Anonymous" rows at the end
#
# Net effect on column A (top to bottom):
#   Helen, Kalli, Chloé, Christina, Carmen, Lee Katherine, Ethan, Emma,
#   Michelle, Gabrielle, Stephanie, Yung, Alex, Anonymous, Anonymous,
#   Anonymous, Anonymous

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Helen",
    "Kalli",
    "Chloé",
    "Christina",
    "Carmen",
    "Lee Katherine",
    "Ethan",
    "Emma",
    "Michelle",
    "Gabrielle",
    "Stephanie",
    "Yung",
    "Alex",
    "Anonymous",
    "Anonymous",
    "Anonymous",
    "Anonymous"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $names[$i]
}
